$wb = $excel.ActiveWorkbook

# ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2324.0715
$ws.Range("I8").Value = 230.63637
$ws.Range("J8").Value = 10000
$ws.Range("K8").Value = 691.9091100000001
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = -552.9091100000001
$ws.Range("N8").Value = -30278

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4448.057
$ws.Range("I137").Value = 986.94116
$ws.Range("J137").Value = 7716.8887
$ws.Range("K137").Value = 2960.82348
$ws.Range("L137").Value = 23150.6661
$ws.Range("M137").Value = -410.82348
$ws.Range("N137").Value = -28250.6661

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4832.0605
$ws.Range("I141").Value = 1736.3103
$ws.Range("J141").Value = 27276.25
$ws.Range("K141").Value = 5208.9309
$ws.Range("L141").Value = 81828.75
$ws.Range("M141").Value = -28.93090000000029
$ws.Range("N141").Value = -92188.75

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 775.7826
$ws.Range("I2").Value = 629.381
$ws.Range("K2").Value = 629.381
$ws.Range("M2").Value = -516.381

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 882.38
$ws.Range("I32").Value = 737.41974
$ws.Range("J32").Value = 1500.3684
$ws.Range("K32").Value = 737.41974
$ws.Range("L32").Value = 1500.3684
$ws.Range("M32").Value = -450.41974
$ws.Range("N32").Value = -2074.3684

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1591.3438
$ws.Range("I97").Value = 640.1739
$ws.Range("J97").Value = 4022.111
$ws.Range("K97").Value = 640.1739
$ws.Range("L97").Value = 4022.111
$ws.Range("M97").Value = -144.1739
$ws.Range("N97").Value = -5014.111

# ARM row 112
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 18500
$ws.Range("J112").Value = 18500
$ws.Range("L112").Value = 18500
$ws.Range("N112").Value = -21454

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 775.7826
$ws.Range("I116").Value = 629.381
$ws.Range("K116").Value = 629.381
$ws.Range("M116").Value = 1664.619

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1924.0952
$ws.Range("I122").Value = 1425.4839
$ws.Range("J122").Value = 3329.2727
$ws.Range("K122").Value = 4276.4517
$ws.Range("L122").Value = 9987.8181
$ws.Range("M122").Value = -1826.4517
$ws.Range("N122").Value = -14887.8181

# ARM row 124
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 36143
$ws.Range("J124").Value = 36143
$ws.Range("L124").Value = 36143
$ws.Range("N124").Value = -45963

# ARM row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 38143
$ws.Range("J125").Value = 38143
$ws.Range("L125").Value = 38143
$ws.Range("N125").Value = -47983

# ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 49400
$ws.Range("J135").Value = 49400
$ws.Range("L135").Value = 49400
$ws.Range("N135").Value = -59540

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 567498.9
$ws.Range("I139").Value = 777000
$ws.Range("J139").Value = 288164
$ws.Range("K139").Value = 777000
$ws.Range("L139").Value = 288164
$ws.Range("M139").Value = -771860
$ws.Range("N139").Value = -298444

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 775.7826
$ws.Range("I3").Value = 629.381
$ws.Range("K3").Value = 629.381
$ws.Range("M3").Value = -515.381

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1397.1904
$ws.Range("I107").Value = 1100.2
$ws.Range("J107").Value = 1667.1818
$ws.Range("K107").Value = 1100.2
$ws.Range("L107").Value = 1667.1818
$ws.Range("M107").Value = 819.8
$ws.Range("N107").Value = -5507.1818

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 47.291668
$ws.Range("J12").Value = 53.38095
$ws.Range("L12").Value = 160.14285
$ws.Range("N12").Value = -506.14285

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1390.9656
$ws.Range("I68").Value = 859.1818
$ws.Range("J68").Value = 1570.9539
$ws.Range("K68").Value = 2577.5454
$ws.Range("L68").Value = 4712.861699999999
$ws.Range("M68").Value = -1766.5454
$ws.Range("N68").Value = -6334.861699999999

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1390.9656
$ws.Range("I71").Value = 859.1818
$ws.Range("J71").Value = 1570.9539
$ws.Range("K71").Value = 7732.6362
$ws.Range("L71").Value = 14138.5851
$ws.Range("M71").Value = -3676.6362
$ws.Range("N71").Value = -22250.5851

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 655.26086
$ws.Range("I107").Value = 405.43103
$ws.Range("J107").Value = 1972.5454
$ws.Range("K107").Value = 1216.29309
$ws.Range("L107").Value = 5917.6362
$ws.Range("M107").Value = 703.7069099999999
$ws.Range("N107").Value = -9757.636200000001

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 591.86
$ws.Range("J131").Value = 914.84
$ws.Range("L131").Value = 2744.52
$ws.Range("N131").Value = -12824.52

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 46136.332
$ws.Range("I11").Value = 22268.666
$ws.Range("J11").Value = 70004
$ws.Range("K11").Value = 22268.666
$ws.Range("L11").Value = 70004
$ws.Range("M11").Value = -22129.666
$ws.Range("N11").Value = -70282

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4806.722
$ws.Range("I113").Value = 5380.1113
$ws.Range("J113").Value = 4233.3335
$ws.Range("K113").Value = 5380.1113
$ws.Range("L113").Value = 4233.3335
$ws.Range("M113").Value = -3210.1113
$ws.Range("N113").Value = -8573.333500000001

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1478
$ws.Range("I7").Value = 1478
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1478
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1366
$ws.Range("N7").ClearContents()

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1385.5
$ws.Range("I61").Value = 1017.3333
$ws.Range("J61").Value = 2490
$ws.Range("K61").Value = 1017.3333
$ws.Range("L61").Value = 2490
$ws.Range("M61").Value = -815.3333
$ws.Range("N61").Value = -2894

# LTW row 110
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 38666.668
$ws.Range("J110").Value = 38666.668
$ws.Range("L110").Value = 38666.668
$ws.Range("N110").Value = -46846.668

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1385.5
$ws.Range("I113").Value = 1017.3333
$ws.Range("J113").Value = 2490
$ws.Range("K113").Value = 1017.3333
$ws.Range("L113").Value = 2490
$ws.Range("M113").Value = 1152.6667
$ws.Range("N113").Value = -6830

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1478
$ws.Range("I126").Value = 1478
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4434
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1964
$ws.Range("N126").ClearContents()

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5738.5
$ws.Range("I132").Value = 6891.4287
$ws.Range("K132").Value = 20674.2861
$ws.Range("M132").Value = -18144.2861
